{"js": "// The sentence \"Dar soluci\u00f3n a los problemas que requieren una respuesta\n// r\u00e1pida en un ambiente flexible y con cambios constantes, haciendo caso\n// omiso de la documentaci\u00f3n rigurosa y los m\u00e9todos formales.\" is left\n// verbatim - only the hidden \"_GoBack\" bookmark (Word's \"last edit\n// location\" marker) needs to move from right before \"r soluci\u00f3n...\" to\n// right after \"...ambiente flexible y\" (i.e. just before \" con cambios\n// constantes...\").\n//\n// Re-inserting the same text over a range that spans the bookmark's\n// current (collapsed) position causes Word to drop the stale mark and\n// re-anchor it at the end of the freshly (re)written text, so we reuse\n// that to relocate \"_GoBack\" without touching any visible content.\n\nconst target =\n  \"ar soluci\u00f3n a los problemas que requieren una respuesta r\u00e1pida en un \" +\n  \"ambiente flexible y\";\n\nconst results = context.document.body.search(target, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the target sentence fragment.\");\n}\n\nconst matchRange = results.items[0];\nmatchRange.insertText(matchRange.text, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The sentence \"Dar soluci\u00f3n a los problemas que requieren una respuesta\n# r\u00e1pida en un ambiente flexible y con cambios constantes, haciendo caso\n# omiso de la documentaci\u00f3n rigurosa y los m\u00e9todos formales.\" keeps its\n# wording - only the hidden \"_GoBack\" bookmark (Word's \"last edit\n# location\" marker) needs to move from right before \"r soluci\u00f3n...\" to\n# right after \"...ambiente flexible y\" (i.e. just before \" con cambios\n# constantes...\").\n\n$d = $word.ActiveDocument\n\n# Locate the existing \"_GoBack\" bookmark and force it to drop: rewriting\n# the two characters that straddle its (collapsed) position with\n# different text, then restoring the original text, removes the stale\n# bookmark without changing any visible content.\n$bm = $d.Bookmarks(\"_GoBack\")\n$bmStart = $bm.Start\n$straddle = $d.Range($bmStart - 1, $bmStart + 1)\n$originalText = $straddle.Text\n$straddle.Text = \"##\"\n$restore = $d.Range($bmStart - 1, $bmStart + 1)\n$restore.Text = $originalText\n\n# Find the new anchor point: right after \"...ambiente flexible y\".\n$findRange = $d.Content\n$findRange.Find.Text = \"ambiente flexible y\"\n$findRange.Find.Execute() | Out-Null\n$newPos = $findRange.End\n\n# Re-create \"_GoBack\" at the new, collapsed location.\n$newBookmarkRange = $d.Range($newPos, $newPos)\n$d.Bookmarks.Add(\"_GoBack\", $newBookmarkRange)\n"}
